$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knowledge")

$ws.Range("E6").Value = "Finance,Data"
$ws.Range("E8").Value = "Finance,Debt,Government"
$ws.Range("E9").Value = "Programming,Python"
$ws.Range("E13").Value = "Finance,Quanitative"
$ws.Range("E14").Value = "Finance,Government"
$ws.Range("E16").Value = "Finance,Blog"
$ws.Range("E17").Value = "Finance,Blog"
$ws.Range("E26").Value = "Programming,Data"
$ws.Range("E27").Value = "Programming,C#,Web Development"
$ws.Range("E28").Value = "Programming,C#,Web Development"
$ws.Range("E29").Value = "Programming, C#"
$ws.Range("E30").Value = "Programming, Web Development"
$ws.Range("E31").Value = "Programming,Data,Python"
$ws.Range("E32").Value = "HTML,CSS, Web Development"
$ws.Range("E34").Value = "C#,Programming,Excel,EPPlus"
$ws.Range("E35").Value = "Python, Programming, Excel"
$ws.Range("E36").Value = "Programming, Database, MongoDB"
$ws.Range("E37").Value = "C#,Programming,Excel,EPPlus"
$ws.Range("E38").Value = "C#,Programming,Python,LINQ"
$ws.Range("E39").Value = "Python, Programming, Pygame,"
$ws.Range("E40").Value = "Python, Programming, FRED"
$ws.Range("E41").Value = "HTML,Web Development"
$ws.Range("E42").Value = "CSS,Web Development"
$ws.Range("E43").Value = "Python, Programming"
$ws.Range("E44").Value = "Programming,Python"
$ws.Range("E49").Value = "Python, Regex"
$ws.Range("E50").Value = "NLP, Language"
$ws.Range("E51").Value = "Python, Programming"
$ws.Range("E52").Value = "Python, Data, Visualization"
$ws.Range("E54").Value = "Python, programming"
$ws.Range("E55").Value = "Finance, blog, research"

# Restore the view state: scroll so A49 is the top-left visible cell,
# and select C57 as the active cell.
$ws.Activate()
$ws.Range("C57").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
